# Batch size experiment analysis
# Append additional batch-size experiment rows (train_size 500 and 300)
# to the results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(500, 1,  0.337386018237082,  0, "0:09:10"),
    @(500, 4,  0.3399735878823331, 4, "0:34:37"),
    @(500, 8,  0.8538578986332499, 6, "0:22:46"),
    @(500, 16, 0.8690897597977244, 9, "0:16:53"),
    @(500, 32, 0.8583499539481968, 9, "0:09:29"),
    @(300, 1,  0.337386018237082,  1, "0:14:34"),
    @(300, 4,  0.8236784407325322, 6, "0:44:18"),
    @(300, 8,  0.8551013625789299, 7, "0:23:21"),
    @(300, 16, 0.8222474726552276, 9, "0:15:10"),
    @(300, 32, 0.8124270846989434, 9, "0:08:44")
)

$startRow = 27
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
